$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get swapped between row 3 and row 4
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr3 = "$col" + "3"
    $addr4 = "$col" + "4"
    $v3 = $ws.Range($addr3).Value2
    $v4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value2 = $v4
    $ws.Range($addr4).Value2 = $v3
}
